$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13, pushing existing rows 13-23 down to 14-24.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly data point.
$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(13, 3).Value = "Maule"
$ws.Cells.Item(13, 4).Value = 45097
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 5).Value = 7
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100104
$ws.Cells.Item(13, 8).Value = "Frutos de pepita"
$ws.Cells.Item(13, 9).Value = 100104001
$ws.Cells.Item(13, 10).Value = "Granada"
$ws.Cells.Item(13, 11).Value = "Wonderfull"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 17000
$ws.Cells.Item(13, 15).Value = 17000
$ws.Cells.Item(13, 16).Value = 17000
$ws.Cells.Item(13, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(13, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 19).Value = 944
$ws.Cells.Item(13, 20).Value = 18
